# Sync updated BOM for multitarget
#
# - The old "Digikey Upload" sheet is retired/renamed to "DO NOT USE".
# - The "Digikey Final" sheet is moved to the front of the workbook.
# - The "From Eagle" sheet becomes the active (selected) tab.
#
# Excel auto-repoints any formulas that reference a renamed sheet, so the
# "Digikey Final" BOM formulas ('Digikey Upload'!H3, etc.) automatically
# become ('DO NOT USE'!H3, etc.) once the rename below runs.

$wb = $excel.ActiveWorkbook

# Rename "Digikey Upload" -> "DO NOT USE" (formulas referencing it update automatically)
$wb.Worksheets.Item("Digikey Upload").Name = "DO NOT USE"

# Reorder: move "Digikey Final" to be the first sheet in the workbook
$wb.Worksheets.Item("Digikey Final").Move($wb.Worksheets.Item(1))

# Make "From Eagle" the active/selected sheet
$wb.Worksheets.Item("From Eagle").Activate()
